$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card9")

# Row 24 currently has empty placeholder cells in B:K and M.
# Fill them with the literal text "nan" to match the rest of the table.
$ws.Range("B24:K24").Value = "nan"
$ws.Range("M24").Value = "nan"

# Add the new event row (25) describing the gearbox oil change.
# Force text storage (matches the "card" column elsewhere, which is text, not numeric),
# then drop back to the default cell style so no stray number format lingers.
$ws.Range("A25").NumberFormat = "@"
$ws.Range("A25").Value = "9"
$ws.Range("A25").Style = "Normal"
$ws.Range("L25").Value = "10\7\2025"
$ws.Range("N25").Value = "تم تغيير زيت الجيربوكس"
$ws.Range("O25").Value = "تيم العمل"
